# Updates cell values across multiple sheets to reflect refreshed
# market-price-derived profit calculations (scheduled data refresh).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 11666.667
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9532
# Row 23
$ws.Range("H23").Value = 11666.667
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9766
# Row 28
$ws.Range("H28").Value = 1017.2308
$ws.Range("I28").Value = 472.9
$ws.Range("J28").Value = 2831.6667
$ws.Range("K28").Value = 472.9
$ws.Range("L28").Value = 2831.6667
$ws.Range("M28").Value = 12.10000000000002
$ws.Range("N28").Value = -3801.6667
# Row 80
$ws.Range("H80").Value = 46579.09
$ws.Range("I80").Value = 1094.2727
$ws.Range("J80").Value = 92063.91
$ws.Range("K80").Value = 3282.8181
$ws.Range("L80").Value = 276191.73
$ws.Range("M80").Value = -2284.8181
$ws.Range("N80").Value = -278187.73
# Row 83
$ws.Range("H83").Value = 46579.09
$ws.Range("I83").Value = 1094.2727
$ws.Range("J83").Value = 92063.91
$ws.Range("K83").Value = 9848.454299999999
$ws.Range("L83").Value = 828575.1900000001
$ws.Range("M83").Value = -4856.454299999999
$ws.Range("N83").Value = -838559.1900000001
# Row 118
$ws.Range("H118").Value = 12121.111
$ws.Range("I118").Value = 17431.666
$ws.Range("J118").Value = 1500
$ws.Range("K118").Value = 52294.99800000001
$ws.Range("L118").Value = 4500
$ws.Range("M118").Value = -50637.99800000001
$ws.Range("N118").Value = -7814
# Row 132
$ws.Range("H132").Value = 5441155
$ws.Range("I132").Value = 5688253
$ws.Range("K132").Value = 17064759
$ws.Range("M132").Value = -17062229
# Row 137
$ws.Range("H137").Value = 959.8611
$ws.Range("I137").Value = 928.1875
$ws.Range("J137").Value = 1213.25
$ws.Range("K137").Value = 2784.5625
$ws.Range("L137").Value = 3639.75
$ws.Range("M137").Value = -234.5625
$ws.Range("N137").Value = -8739.75
# Row 140
$ws.Range("H140").Value = 64700
$ws.Range("J140").Value = 64700
$ws.Range("L140").Value = 64700
$ws.Range("N140").Value = -75060
# Row 141
$ws.Range("H141").Value = 3824.2856
$ws.Range("I141").Value = 3826.3333
$ws.Range("J141").Value = 3822.75
$ws.Range("K141").Value = 11478.9999
$ws.Range("L141").Value = 11468.25
$ws.Range("M141").Value = -6298.999899999999
$ws.Range("N141").Value = -21828.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 3999.6667
$ws.Range("J8").Value = 3999.6667
$ws.Range("L8").Value = 3999.6667
$ws.Range("N8").Value = -4287.6667
# Row 32
$ws.Range("H32").Value = 23980.09
$ws.Range("I32").Value = 6711.284
$ws.Range("K32").Value = 6711.284
$ws.Range("M32").Value = -6424.284
# Row 68
$ws.Range("H68").Value = 44990
$ws.Range("J68").Value = 44990
$ws.Range("L68").Value = 44990
$ws.Range("N68").Value = -46612
# Row 71
$ws.Range("H71").Value = 44990
$ws.Range("J71").Value = 44990
$ws.Range("L71").Value = 134970
$ws.Range("N71").Value = -143082
# Row 96
$ws.Range("H96").Value = 21433.334
$ws.Range("J96").Value = 21433.334
$ws.Range("L96").Value = 21433.334
$ws.Range("N96").Value = -26925.334
# Row 132
$ws.Range("H132").Value = 14560.319
$ws.Range("I132").Value = 18474.53
$ws.Range("K132").Value = 55423.59
$ws.Range("M132").Value = -52893.59

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 93759.586
$ws.Range("I86").Value = 139605
$ws.Range("J86").Value = 2068.75
$ws.Range("K86").Value = 139605
$ws.Range("L86").Value = 2068.75
$ws.Range("M86").Value = -138482
$ws.Range("N86").Value = -4314.75
# Row 89
$ws.Range("H89").Value = 93759.586
$ws.Range("I89").Value = 139605
$ws.Range("J89").Value = 2068.75
$ws.Range("K89").Value = 698025
$ws.Range("L89").Value = 10343.75
$ws.Range("M89").Value = -692409
$ws.Range("N89").Value = -21575.75
# Row 128
$ws.Range("H128").Value = 1000040
$ws.Range("I128").Value = 1000040
$ws.Range("K128").Value = 3000120
$ws.Range("M128").Value = -2997630
# Row 134
$ws.Range("H134").Value = 12003.3125
$ws.Range("I134").Value = 13283.536
$ws.Range("J134").Value = 4504.857
$ws.Range("K134").Value = 39850.608
$ws.Range("L134").Value = 13514.571
$ws.Range("M134").Value = -37315.608
$ws.Range("N134").Value = -18584.571

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 63
$ws.Range("H63").Value = 42990
$ws.Range("J63").Value = 42990
$ws.Range("L63").Value = 42990
$ws.Range("N63").Value = -44362
# Row 64
$ws.Range("H64").Value = 40013
$ws.Range("J64").Value = 40013
$ws.Range("L64").Value = 40013
$ws.Range("N64").Value = -40509
# Row 66
$ws.Range("H66").Value = 42990
$ws.Range("J66").Value = 42990
$ws.Range("L66").Value = 128970
$ws.Range("N66").Value = -135834
# Row 67
$ws.Range("H67").Value = 40013
$ws.Range("J67").Value = 40013
$ws.Range("L67").Value = 40013
$ws.Range("N67").Value = -41729
# Row 68
$ws.Range("H68").Value = 17513.047
$ws.Range("J68").Value = 17513.047
$ws.Range("L68").Value = 17513.047
$ws.Range("N68").Value = -19011.047
# Row 69
$ws.Range("H69").Value = 29800
$ws.Range("I69").Value = 29800
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 29800
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -29051
$ws.Range("N69").ClearContents()
# Row 70
$ws.Range("H70").Value = 11999.6
$ws.Range("J70").Value = 11999.6
$ws.Range("L70").Value = 11999.6
$ws.Range("N70").Value = -12629.6
# Row 71
$ws.Range("H71").Value = 17513.047
$ws.Range("J71").Value = 17513.047
$ws.Range("L71").Value = 52539.141
$ws.Range("N71").Value = -60027.141
# Row 72
$ws.Range("H72").Value = 29800
$ws.Range("I72").Value = 29800
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 89400
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -85656
$ws.Range("N72").ClearContents()
# Row 73
$ws.Range("H73").Value = 11999.6
$ws.Range("J73").Value = 11999.6
$ws.Range("L73").Value = 11999.6
$ws.Range("N73").Value = -14183.6
# Row 74
$ws.Range("H74").Value = 32615.555
$ws.Range("J74").Value = 32615.555
$ws.Range("L74").Value = 32615.555
$ws.Range("N74").Value = -34363.555
# Row 77
$ws.Range("H77").Value = 32615.555
$ws.Range("J77").Value = 32615.555
$ws.Range("L77").Value = 97846.66500000001
$ws.Range("N77").Value = -106582.665
# Row 122
$ws.Range("H122").Value = 1237.75
$ws.Range("I122").Value = 1404
$ws.Range("J122").Value = 1138
$ws.Range("K122").Value = 4212
$ws.Range("L122").Value = 3414
$ws.Range("M122").Value = -1762
$ws.Range("N122").Value = -8314

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1197.8096
$ws.Range("I5").Value = 1084.8667
$ws.Range("J5").Value = 1260.5555
$ws.Range("K5").Value = 3254.6001
$ws.Range("L5").Value = 3781.6665
$ws.Range("M5").Value = -3142.6001
$ws.Range("N5").Value = -4005.6665
# Row 129
$ws.Range("H129").Value = 5575.6665
$ws.Range("I129").Value = 622
$ws.Range("J129").Value = 6879.263
$ws.Range("K129").Value = 1866
$ws.Range("L129").Value = 20637.789
$ws.Range("M129").Value = 3134
$ws.Range("N129").Value = -30637.789
# Row 135
$ws.Range("H135").Value = 1197.8096
$ws.Range("I135").Value = 1084.8667
$ws.Range("J135").Value = 1260.5555
$ws.Range("K135").Value = 9763.800300000001
$ws.Range("L135").Value = 11344.9995
$ws.Range("M135").Value = -7228.800300000001
$ws.Range("N135").Value = -16414.9995

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 111625.055
$ws.Range("I70").Value = 186977.81
$ws.Range("K70").Value = 186977.81
$ws.Range("M70").Value = -186707.81
# Row 73
$ws.Range("H73").Value = 111625.055
$ws.Range("I73").Value = 186977.81
$ws.Range("K73").Value = 186977.81
$ws.Range("M73").Value = -186041.81
# Row 102
$ws.Range("H102").Value = 433509
$ws.Range("I102").Value = 4660
$ws.Range("J102").Value = 550467.8
$ws.Range("K102").Value = 4660
$ws.Range("L102").Value = 550467.8
$ws.Range("M102").Value = -3038
$ws.Range("N102").Value = -553711.8
# Row 122
$ws.Range("H122").Value = 3850
$ws.Range("I122").Value = 4466.6665
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 13399.9995
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -10949.9995
$ws.Range("N122").Value = -10900
# Row 126
$ws.Range("H126").Value = 6539154.5
$ws.Range("J126").Value = 11766920
$ws.Range("L126").Value = 35300760
$ws.Range("N126").Value = -35305700
# Row 132
$ws.Range("H132").Value = 2718.0256
$ws.Range("I132").Value = 2147.5334
$ws.Range("J132").Value = 4619.6665
$ws.Range("K132").Value = 6442.600199999999
$ws.Range("L132").Value = 13858.9995
$ws.Range("M132").Value = -3912.600199999999
$ws.Range("N132").Value = -18918.9995

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 237352.11
$ws.Range("J55").Value = 462.8846
$ws.Range("L55").Value = 462.8846
$ws.Range("N55").Value = -808.8846
# Row 61
$ws.Range("H61").Value = 1282.1562
$ws.Range("I61").Value = 1234.9
$ws.Range("J61").Value = 1360.9166
$ws.Range("K61").Value = 1234.9
$ws.Range("L61").Value = 1360.9166
$ws.Range("M61").Value = -1032.9
$ws.Range("N61").Value = -1764.9166
# Row 113
$ws.Range("H113").Value = 1282.1562
$ws.Range("I113").Value = 1234.9
$ws.Range("J113").Value = 1360.9166
$ws.Range("K113").Value = 1234.9
$ws.Range("L113").Value = 1360.9166
$ws.Range("M113").Value = 935.0999999999999
$ws.Range("N113").Value = -5700.9166

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 80011
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 126
$ws.Range("H126").Value = 1084.1538
$ws.Range("I126").Value = 1160.8235
$ws.Range("J126").Value = 939.3333
$ws.Range("K126").Value = 3482.4705
$ws.Range("L126").Value = 2817.9999
$ws.Range("M126").Value = -1012.4705
$ws.Range("N126").Value = -7757.9999
# Row 132
$ws.Range("H132").Value = 3210.361
$ws.Range("I132").Value = 3525.7693
$ws.Range("K132").Value = 10577.3079
$ws.Range("M132").Value = -8047.3079
